$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A134").Value = "2023-12-09 09:44:06"
$ws.Range("B134").Value = 0.0014

$ws.Range("A135").Value = "2023-12-09 09:44:17"
$ws.Range("B135").Value = 0.0004
